$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q: header + per-department case counts for "T15: 2/4/2020".
$ws.Range("Q1").Value = "T15: 2/4/2020"

$ws.Range("Q2").Value = 3     # Atlantida
$ws.Range("Q3").Value = 2     # Choluteca
$ws.Range("Q4").Value = 18    # Colon
$ws.Range("Q5").Value = 0     # Comayagua
$ws.Range("Q6").Value = 0     # Copan
$ws.Range("Q7").Value = 140   # Cortes
$ws.Range("Q8").Value = 0     # El Paraiso
$ws.Range("Q9").Value = 43    # Fco Morazan
$ws.Range("Q10").Value = 0    # Gracias a Dios
$ws.Range("Q11").Value = 0    # Intibuca
$ws.Range("Q12").Value = 0    # Islas de la Bahia
$ws.Range("Q13").Value = 0    # La Paz
$ws.Range("Q14").Value = 4    # Lempira
$ws.Range("Q15").Value = 0    # Ocotepeque
$ws.Range("Q16").Value = 0    # Olancho
$ws.Range("Q17").Value = 4    # Santa Barbara
$ws.Range("Q18").Value = 0    # Valle
$ws.Range("Q19").Value = 5    # Yoro

# Column total, same pattern as the other date columns.
$ws.Range("Q20").Formula = "=SUM(Q2:Q19)"

# Give the new column roughly the same width the author set for it.
$ws.Columns.Item(17).ColumnWidth = 16.5

# Leave the selection where the author left it when they finished editing.
$ws.Range("Q18").Select()
